# Update the 合肥-漫展信息 workbook: several exhibition listings changed.
# The first three events (环形宇宙动漫游戏嘉年华, 环形宇宙动漫游戏嘉年华—吴晛专场,
# 巢湖·原×铁×崩only) are removed, every remaining event shifts up, the
# sequence number column is renumbered, and a handful of "想去人数" (view
# count) figures are refreshed.  Both the "展览" sheet and the "全部类型"
# sheet carry the same table, so the edit is applied to each.

$wb = $excel.ActiveWorkbook
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Remove the three leading rows that no longer appear in the listing.
    # Deleting row 2 three times shifts everything below up automatically.
    $ws.Rows.Item(2).Delete()
    $ws.Rows.Item(2).Delete()
    $ws.Rows.Item(2).Delete()

    # Renumber the sequence column (A2:A10 => 1..9).
    for ($i = 0; $i -lt 9; $i++) {
        $ws.Cells.Item(2 + $i, 1).Value = $i + 1
    }

    # Refresh the "想去人数" (F column) figures that changed for the
    # surviving rows.
    $ws.Cells.Item(3, 6).Value = 2100    # 合肥·新春AG动漫游戏盛典热血plus
    $ws.Cells.Item(4, 6).Value = 1594    # 合肥·2024运动新春动漫庆典（全ip）
    $ws.Cells.Item(5, 6).Value = 316     # 合肥·安徽马娘only
    $ws.Cells.Item(6, 6).Value = 1026    # 合肥·星芒1.5动漫嘉年华
    $ws.Cells.Item(7, 6).Value = 482     # 合肥·CW国潮动漫游戏嘉年华
    $ws.Cells.Item(9, 6).Value = 5660    # 合肥· 第二届漫画城市动漫展 -故事再次开始
}
